$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Objetivos:" value (B10/C10) was mistakenly replaced with the
# "Docentes responsaveis:" value in the source data.
$ws.Range("B10:C10").Value = "5817066 - Félix Monteiro Pereira"

# The old stand-alone row 13 (B13/C13 only, no A label) held the
# "Docentes responsaveis:" value and is removed entirely; everything
# below shifts up by one row.
$ws.Rows.Item(13).Delete()

# After the shift, re-align the remaining content values that no
# longer correspond 1:1 with their (now one-row-higher) labels.
$ws.Range("B13:C13").Value = "Semestral"

$ws.Range("B15:C15").Value = "01/01/2018"

$ws.Range("B18:C18").Value = "5817066 - Félix Monteiro Pereira"

$ws.Range("B19:C19").Value = "A avaliação do aprendizado será realizada por meio da resolução de problemas de engenharia química propostos aos alunos. A média final será calculada pela média aritmética entre a nota obtida na resolução de problemas relativos aos itens 1 e 2 do programa do curso."

$ws.Range("B20:C20").Value = "MF = (P1 + P2)/2 Onde: P1 é a nota obtida pela avaliação da resolução de problemas referentes ao item 1 do Programa do curso;P2 é a nota obtida pela avaliação da resolução de problemas referentes ao item 2 do Programa do curso;MF é a média final do período."

$ws.Range("B21:C21").Value = "A recuperação será feita por meio de uma prova (PR) para alunos que tenham MF maior ou igual a 3,0 e menor do que 5,0 e pelo menos 70% de frequência. A nota de recuperação (NR) será a média simples entre a média final (MF) e a prova de recuperação (PR). Será considerado aprovado o aluno com NR maior ou igual a 5,0."
